$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are designing a database and you want to add a new tuple to a relation that does not contain any foreign key. You want to check all the constraints that might be violated because of that operation. Which constraints should you examine to cover all possible constraint violation cases?",
        "ques_type": 15,
        "options": [
            "Referential integrity constraint",
            "Key constraint",
            "Domain constraint",
            "Entity integrity constraint",
            "Serial constraint"
        ],
        "score": [
            "Key constraint",
            "Domain constraint"
        ]
    },
    {
        "title": "You have the following functional dependencies in R(ABCD):{ A -&gt B, B -&gt C, C -&gt D} Which decomposition of R is going to be lossless?",
        "ques_type": 2,
        "options": [
            "R1(A, B, C) and R2(D)",
            "R1(A, B) and R2(C, D)",
            "R1(A, B) and R2(B, C) and R3(C, D)",
            "R1(A) and R2(B) and R3(C, D)"
        ],
        "score": "R1(A, B) and R2(B, C) and R3(C, D)"
    },
    {
        "title": "You need to apply four transactions related to banking payments. You want to apply read operations to all of them. You are going to execute transactions concurrently. How many transactions may have inconsistent outcomes if scheduling is applied randomly?",
        "ques_type": 2,
        "options": [
            "0",
            "1",
            "3",
            "4"
        ],
        "score": "3"
    },
    {
        "title": "You have designed a B+ tree with an order of three, and in each of its four leaf nodes you have a single index entry. You are going to insert a new index entry into the tree. How many of the leaf nodes might get overflow after trying to apply insertion?",
        "ques_type": 2,
        "options": [
            "0",
            "2",
            "v",
            "4"
        ],
        "score": "0"
    }
]
'@

# Remove trailing newline introduced by the here-string
$text = $text.TrimEnd("`r","`n")

# A2 no longer holds any content
$ws.Range("A2").ClearContents()

# A1 drops its bold/border/centered formatting and becomes the (formatted) questions text
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
